# Q3 Update - 2025
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Global shared-string update: "gtH5M7" -> "abK5LD" for every data row's
#    short-url column (B2:B239 before the row insertion below).
$ws.Range("B2:B239").Value = "abK5LD"

# 2) Insert a brand new data row for Nigeria (coo_id 141) at row 232,
#    which pushes the Rwanda..Yemen rows down by one (232->233 ... 239->240).
$ws.Rows.Item(232).Insert()

# After the insert, row 232 is blank (only formats were copied down).
# Re-assert the short-url value and fill in the new Nigeria record.
$ws.Range("B232").Value = "abK5LD"

$ws.Range("A232").Value = "1"
$ws.Range("C232").Value = "1"
$ws.Range("D232").Value = "231"
$ws.Range("E232").Value = "2024"
$ws.Range("F232").Value = "141"
$ws.Range("G232").Value = "Nigeria"
$ws.Range("H232").Value = "NIG"
$ws.Range("I232").Value = "NGA"
$ws.Range("J232").Value = "68"
$ws.Range("K232").Value = "Gambia"
$ws.Range("L232").Value = "GAM"
$ws.Range("M232").Value = "GMB"
$ws.Range("N232").Value = "0"
$ws.Range("O232").Value = "5"
$ws.Range("P232").Value = "0"
$ws.Range("Q232").Value = "0"
$ws.Range("R232").Value = "0"
$ws.Range("S232").Value = "0"
$ws.Range("T232").Value = "0"
$ws.Range("U232").Value = "-"
$ws.Range("V232").Value = "0"

# 3) Fix up the "items" id column (D) for the rows that shifted down one
#    position (their D value must advance by one to stay row-1).
$ws.Range("D233").Value = "232"
$ws.Range("D234").Value = "233"
$ws.Range("D235").Value = "234"
$ws.Range("D236").Value = "235"
$ws.Range("D237").Value = "236"
$ws.Range("D238").Value = "237"
$ws.Range("D239").Value = "238"
$ws.Range("D240").Value = "239"

# 4) Apply the Q3 2025 refugee/asylum-seeker figure updates.

# Afghanistan (row 225): refugees 0->5, asylum_seekers 5->0
$ws.Range("N225").Value = "5"
$ws.Range("O225").Value = "0"

# Central African Rep. (row 226): refugees 11->13
$ws.Range("N226").Value = "13"

# Congo (row 228): asylum_seekers 6->11
$ws.Range("O228").Value = "11"

# Dem. Rep. of the Congo (row 229): refugees 53->55, asylum_seekers 0->5
$ws.Range("N229").Value = "55"
$ws.Range("O229").Value = "5"

# Liberia (row 231): refugees 15->19
$ws.Range("N231").Value = "19"

# Senegal (now row 234): refugees 3661->3721, asylum_seekers 379->399
$ws.Range("N234").Value = "3721"
$ws.Range("O234").Value = "399"

# Sierra Leone (now row 235): refugees 99->100
$ws.Range("N235").Value = "100"

# Sudan (now row 237): asylum_seekers 12->15
$ws.Range("O237").Value = "15"
